$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEC-2020")

# ---------------------------------------------------------------------
# 1. Stash the formatting of the existing legend block (rows 19-23,
#    columns B & C) into scratch cells on row 1, far outside the table,
#    so it survives the row deletion below untouched.
# ---------------------------------------------------------------------
$legendRows = 19, 20, 21, 22, 23
$legendText = @{
    "19,3" = "Status"
    "20,2" = "WIP"; "20,3" = "WIP"
    "21,2" = "Pending"; "21,3" = "Pending"
    "22,2" = "Completed"; "22,3" = "Completed"
    "23,2" = "Hold"; "23,3" = "Hold"
}
$scratchCol = 30
foreach ($r in $legendRows) {
    foreach ($c in 2, 3) {
        $key = "$r,$c"
        if ($legendText.ContainsKey($key)) {
            $ws.Cells.Item($r, $c).Copy()
            $ws.Cells.Item(1, $scratchCol).PasteSpecial(-4122)
            $scratchCol = $scratchCol + 1
        }
    }
}

# ---------------------------------------------------------------------
# 2. Remove the old legend rows outright (nothing below them, so the
#    shift collapses to a clean delete with no leftover row nodes).
# ---------------------------------------------------------------------
$ws.Rows("19:23").Delete()

# ---------------------------------------------------------------------
# 3. Append three more tracker entries (rows 17-19), matching the
#    formatting already used by row 16 directly above, with the Status
#    column styled like the existing "Completed" entries (row 2).
# ---------------------------------------------------------------------
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G19").PasteSpecial(-4122)
$ws.Cells.Item(2, 6).Copy()
$ws.Range("F17:F19").PasteSpecial(-4122)

$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = 44181
$ws.Cells.Item(17, 3).Value = "QMVAR 2.0"
$ws.Cells.Item(17, 4).Value = "Design issue fixing"
$ws.Cells.Item(17, 6).Value = "Completed"

$ws.Cells.Item(18, 1).Value = 8
$ws.Cells.Item(18, 2).Value = 44182
$ws.Cells.Item(18, 3).Value = "QMVAR 2.0"
$ws.Cells.Item(18, 4).Value = "Design issue fixing in scroll bar"
$ws.Cells.Item(18, 6).Value = "Completed"

$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = 44183
$ws.Cells.Item(19, 3).Value = "QMVAR 2.0"
$ws.Cells.Item(19, 4).Value = "Design issue fixing in dialog box"
$ws.Cells.Item(19, 6).Value = "Completed"

# ---------------------------------------------------------------------
# 4. Re-create the legend 8 rows further down (27-31), pulling the
#    stashed formatting back in, then writing its text back in place.
# ---------------------------------------------------------------------
$scratchCol = 30
foreach ($r in $legendRows) {
    $dstRow = $r + 8
    foreach ($c in 2, 3) {
        $key = "$r,$c"
        if ($legendText.ContainsKey($key)) {
            $ws.Cells.Item(1, $scratchCol).Copy()
            $ws.Cells.Item($dstRow, $c).PasteSpecial(-4122)
            $ws.Cells.Item($dstRow, $c).Value = $legendText[$key]
            $scratchCol = $scratchCol + 1
        }
    }
}

# Clean up the scratch cells used to carry formatting across the delete.
$ws.Range($ws.Cells.Item(1, 30), $ws.Cells.Item(1, $scratchCol - 1)).Clear()

# ---------------------------------------------------------------------
# 5. Match the view state Excel recorded after the edit: scrolled so
#    row 9 is at the top, with the active cell on D20.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("D20").Select()
